$wb = $excel.ActiveWorkbook

$wsUi = $wb.Worksheets.Item("c-demo_ui")
$wsTest = $wb.Worksheets.Item("t-testCase1")

# Update the "omdb call" row on the test-case sheet to call the new named
# path-based request instead of the old literal endpoint/verb pair.
$wsTest.Range("C3").Value = "omdb.Wall-E"
$wsTest.Range("D3").Value = ""
$wsTest.Range("E3").Value = "export::wall-e_plot::JSONPATH::Plot"
$wsTest.Range("F3").Value = ""

# Restore the selection/view on the test-case sheet and make the config
# sheet the active tab again, matching the saved workbook state.
$wsTest.Range("F3").Select() | Out-Null
$wsUi.Activate() | Out-Null
